$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to stage new text values so that numeric-looking
# strings (e.g. "244.27") are pasted back in as TEXT, matching the
# original inline-string cell type, instead of being auto-converted
# to numbers by the normal .Value assignment parser.
$helper = $ws.Range("H100")
$helper.NumberFormat = "@"

function Set-TextValue($cell, $text) {
    $helper.Value = $text
    $helper.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue $ws.Cells.Item(2, 4) "36.670.49"
Set-TextValue $ws.Cells.Item(2, 5) "  +0.92%  "

Set-TextValue $ws.Cells.Item(3, 4) "1.961.35"
Set-TextValue $ws.Cells.Item(3, 5) "  +1.57%  "

Set-TextValue $ws.Cells.Item(4, 5) "  -0.08%  "

Set-TextValue $ws.Cells.Item(5, 4) "244.27"
Set-TextValue $ws.Cells.Item(5, 5) "  +1.52%  "

Set-TextValue $ws.Cells.Item(6, 4) "0.616"
Set-TextValue $ws.Cells.Item(6, 5) "  +1.74%  "

Set-TextValue $ws.Cells.Item(7, 4) "60.25"
Set-TextValue $ws.Cells.Item(7, 5) "  +7.48%  "

Set-TextValue $ws.Cells.Item(8, 5) "  +0.00%  "

Set-TextValue $ws.Cells.Item(9, 4) "0.377"
Set-TextValue $ws.Cells.Item(9, 5) "  +5.82%  "

Set-TextValue $ws.Cells.Item(10, 4) "0.0792"
Set-TextValue $ws.Cells.Item(10, 5) "  -4.61%  "

Set-TextValue $ws.Cells.Item(11, 5) "  -0.07%  "

Set-TextValue $ws.Cells.Item(12, 4) "14.22"
Set-TextValue $ws.Cells.Item(12, 5) "  +7.40%  "

Set-TextValue $ws.Cells.Item(13, 4) "0.842"
Set-TextValue $ws.Cells.Item(13, 5) "  +5.78%  "

Set-TextValue $ws.Cells.Item(14, 4) "2.248.72"
Set-TextValue $ws.Cells.Item(14, 5) "  +1.53%  "

Set-TextValue $ws.Cells.Item(15, 4) "21.48"
Set-TextValue $ws.Cells.Item(15, 5) "  +3.85%  "

Set-TextValue $ws.Cells.Item(16, 4) "5.31"
Set-TextValue $ws.Cells.Item(16, 5) "  +4.19%  "

Set-TextValue $ws.Cells.Item(17, 4) "1.962.59"
Set-TextValue $ws.Cells.Item(17, 5) "  +1.42%  "

Set-TextValue $ws.Cells.Item(18, 4) "36.587.90"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.91%  "

Set-TextValue $ws.Cells.Item(19, 4) "69.79"
Set-TextValue $ws.Cells.Item(19, 5) "  +1.72%  "

Set-TextValue $ws.Cells.Item(20, 4) "0.0₃0852"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.42%  "

Set-TextValue $ws.Cells.Item(21, 4) "230.14"
Set-TextValue $ws.Cells.Item(21, 5) "  +1.88%  "

Set-TextValue $ws.Cells.Item(22, 4) "5.09"
Set-TextValue $ws.Cells.Item(22, 5) "  +3.62%  "

Set-TextValue $ws.Cells.Item(23, 5) "  +0.10%  "

Set-TextValue $ws.Cells.Item(24, 4) "2.47"
Set-TextValue $ws.Cells.Item(24, 5) "  +6.41%  "

Set-TextValue $ws.Cells.Item(25, 5) "  +4.85%  "

Set-TextValue $ws.Cells.Item(26, 4) "0.146"
Set-TextValue $ws.Cells.Item(26, 5) "  +12.40%  "

Set-TextValue $ws.Cells.Item(27, 4) "9.22"
Set-TextValue $ws.Cells.Item(27, 5) "  +1.53%  "

Set-TextValue $ws.Cells.Item(28, 4) "160.72"
Set-TextValue $ws.Cells.Item(28, 5) "  +0.30%  "

Set-TextValue $ws.Cells.Item(29, 4) "19.50"
Set-TextValue $ws.Cells.Item(29, 5) "  +2.37%  "

Set-TextValue $ws.Cells.Item(30, 4) "1.21"
Set-TextValue $ws.Cells.Item(30, 5) "  +10.77%  "

Set-TextValue $ws.Cells.Item(32, 4) "4.76"
Set-TextValue $ws.Cells.Item(32, 5) "  +5.45%  "

Set-TextValue $ws.Cells.Item(33, 4) "0.0618"
Set-TextValue $ws.Cells.Item(33, 5) "  -0.52%  "

Set-TextValue $ws.Cells.Item(34, 5) "  +7.74%  "

Set-TextValue $ws.Cells.Item(35, 4) "3.58"
Set-TextValue $ws.Cells.Item(35, 5) "  +21.60%  "

Set-TextValue $ws.Cells.Item(36, 4) "2.29"
Set-TextValue $ws.Cells.Item(36, 5) "  +8.38%  "

Set-TextValue $ws.Cells.Item(37, 4) "0.999"
Set-TextValue $ws.Cells.Item(37, 5) "  -0.16%  "

Set-TextValue $ws.Cells.Item(38, 5) "  -1.05%  "

Set-TextValue $ws.Cells.Item(39, 4) "5.55"
Set-TextValue $ws.Cells.Item(39, 5) "  -7.46%  "

Set-TextValue $ws.Cells.Item(40, 4) "0.0983"
Set-TextValue $ws.Cells.Item(40, 5) "  +2.09%  "

Set-TextValue $ws.Cells.Item(41, 5) "  +1.34%  "

Set-TextValue $ws.Cells.Item(42, 5) "  +3.65%  "

Set-TextValue $ws.Cells.Item(43, 5) "  +1.84%  "

Set-TextValue $ws.Cells.Item(44, 4) "16.03"
Set-TextValue $ws.Cells.Item(44, 5) "  +4.35%  "

Set-TextValue $ws.Cells.Item(45, 4) "1.371.62"
Set-TextValue $ws.Cells.Item(45, 5) "  +3.46%  "

Set-TextValue $ws.Cells.Item(48, 4) "7.15"
Set-TextValue $ws.Cells.Item(48, 5) "  +1.82%  "

Set-TextValue $ws.Cells.Item(49, 4) "2.85"
Set-TextValue $ws.Cells.Item(49, 5) "  +0.83%  "

Set-TextValue $ws.Cells.Item(50, 4) "44.64"
Set-TextValue $ws.Cells.Item(50, 5) "  +3.34%  "

Set-TextValue $ws.Cells.Item(51, 4) "2.138.52"
Set-TextValue $ws.Cells.Item(51, 5) "  +1.50%  "

# Row 46 / 47: Aave and ARBITRUM swap positions (with refreshed price/volume)
Set-TextValue $ws.Cells.Item(46, 2) "ARBITRUM"
Set-TextValue $ws.Cells.Item(46, 3) "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Cells.Item(46, 4) "1.03"
Set-TextValue $ws.Cells.Item(46, 5) "  +2.54%  "
Set-TextValue $ws.Cells.Item(47, 2) "Aave"
Set-TextValue $ws.Cells.Item(47, 3) "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(47, 4) "88.38"
Set-TextValue $ws.Cells.Item(47, 5) "  +4.17%  "

$helper.Clear()

Write-Host "Updated cryptos list"